# Extend the cache-conflict benchmark ("cc_results") data to cover
# larger working sets (2^11..2^13), refresh the already-run sizes with
# the new measurements, and re-order the "aligned" block so the new
# 2^11..2^13 "unaligned" rows sit right after the original 11
# "unaligned" rows, followed by the full (now 14-row) "aligned" block.
#
# cc_charts and cc_plot (and the line chart built on cc_plot) are all
# formula-driven off of cc_results, so once the raw numbers here are
# updated Excel's own recalculation repopulates everything downstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cc_results")

# Timestamp in A2 moved forward (benchmark re-run).
$ws.Range("A2").Value = 43028.162870370368

# Columns A,B,C,D,E,K,L for rows 4..31.
$data = New-Object 'object[,]' 28,7

function Set-Row([int]$i, [string]$name, [double]$b, [double]$c, [double]$d, [int]$k, [int]$l) {
    $data[$i,0] = $name
    $data[$i,1] = $b
    $data[$i,2] = $c
    $data[$i,3] = $d
    $data[$i,4] = "ns"
    $data[$i,5] = $k
    $data[$i,6] = $l
}

# "unaligned" (working_set / K=1) rows, now 0..13 (was 0..10).
Set-Row 0  "BM_cache_conflict/0/1"  497777778   1.3259799999999999  1.31836              1 0
Set-Row 1  "BM_cache_conflict/1/1"  213333333   3.2218100000000001  3.2226599999999999  1 1
Set-Row 2  "BM_cache_conflict/2/1"  112000000   6.1873300000000002  6.1383900000000002  1 2
Set-Row 3  "BM_cache_conflict/3/1"  56000000    12.3729              12.2768             1 3
Set-Row 4  "BM_cache_conflict/4/1"  24888889    29.255299999999998   28.878299999999999 1 4
Set-Row 5  "BM_cache_conflict/5/1"  10000000    52.867400000000004   53.125              1 5
Set-Row 6  "BM_cache_conflict/6/1"  5600000     119.119               119.97799999999999 1 6
Set-Row 7  "BM_cache_conflict/7/1"  1544828     455.32299999999998   455.14800000000002 1 7
Set-Row 8  "BM_cache_conflict/8/1"  746667      936.45799999999997   941.68499999999995 1 8
Set-Row 9  "BM_cache_conflict/9/1"  373333      1922.39              1925.22             1 9
Set-Row 10 "BM_cache_conflict/10/1" 100000      5397.44              5312.5              1 10
Set-Row 11 "BM_cache_conflict/11/1" 40727       16368.6              16497               1 11
Set-Row 12 "BM_cache_conflict/12/1" 10000       55459.9              54687.5             1 12
Set-Row 13 "BM_cache_conflict/13/1" 6400        110953               112305              1 13

# "aligned" (working_set / K=0) rows, now 0..13 (was 0..10).
Set-Row 14 "BM_cache_conflict/0/0"  560000000   1.32623              1.31138             0 0
Set-Row 15 "BM_cache_conflict/1/0"  224000000   3.2066699999999999  3.20871             0 1
Set-Row 16 "BM_cache_conflict/2/0"  112000000   6.1869399999999999  6.2778999999999998 0 2
Set-Row 17 "BM_cache_conflict/3/0"  56000000    12.372400000000001   12.5558             0 3
Set-Row 18 "BM_cache_conflict/4/0"  19478261    36.692100000000003   36.900100000000002 0 4
Set-Row 19 "BM_cache_conflict/5/0"  8960000     85.726399999999998   85.449200000000005 0 5
Set-Row 20 "BM_cache_conflict/6/0"  3446154     195.74700000000001   194.964             0 6
Set-Row 21 "BM_cache_conflict/7/0"  560000      1313.22              1311.38             0 7
Set-Row 22 "BM_cache_conflict/8/0"  248889      2797.74              2825.05             0 8
Set-Row 23 "BM_cache_conflict/9/0"  100000      5698.97              5625                0 9
Set-Row 24 "BM_cache_conflict/10/0" 64000       12011.7              11962.9             0 10
Set-Row 25 "BM_cache_conflict/11/0" 10000       49537.1              50000               0 11
Set-Row 26 "BM_cache_conflict/12/0" 3200        235429               234375              0 12
Set-Row 27 "BM_cache_conflict/13/0" 1000        505277               515625              0 13

# Columns F..J are unused for the data block (only the header row uses
# them); E/K/L columns flank that unused gap, so write the two used
# column groups (A:E and K:L) separately to avoid touching F:J.
$abcde = New-Object 'object[,]' 28,5
$kl = New-Object 'object[,]' 28,2
for ($i = 0; $i -lt 28; $i++) {
    $abcde[$i,0] = $data[$i,0]
    $abcde[$i,1] = $data[$i,1]
    $abcde[$i,2] = $data[$i,2]
    $abcde[$i,3] = $data[$i,3]
    $abcde[$i,4] = $data[$i,4]
    $kl[$i,0] = $data[$i,5]
    $kl[$i,1] = $data[$i,6]
}
$ws.Range("A4:E31").Value = $abcde
$ws.Range("K4:L31").Value = $kl

# View state: whole new data block selected, top-left reset.
$ws.Range("A1:L31").Select()
